$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) values
$ws.Cells.Item(2, 4).Value = "61.948.46"
$ws.Cells.Item(3, 4).Value = "3.415.84"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "575.90"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "148.67"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "8.05"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.413"
$ws.Cells.Item(12, 4).Value = "3.998.57"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "28.32"
$ws.Cells.Item(15, 4).Value = "3.404.48"
$ws.Cells.Item(17, 4).Value = "61.930.00"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "14.40"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "8.87"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "380.13"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "75.13"
$ws.Cells.Item(25, 4).Value = "3.555.74"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "7.63"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.998"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "7.89"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.999"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "5.46"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "30.86"
$ws.Cells.Item(40, 4).Value = "3.447.97"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0774"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "42.68"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.775"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.16"
$ws.Cells.Item(47, 4).Value = "2.541.10"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "6.87"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "22.37"

# Update Volume(1h) column (E) values
$ws.Cells.Item(2, 5).Value = "  -2.15%  "
$ws.Cells.Item(3, 5).Value = "  -1.61%  "
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
$ws.Cells.Item(5, 5).Value = "  -0.73%  "
$ws.Cells.Item(6, 5).Value = "  +0.70%  "
$ws.Cells.Item(7, 5).Value = "  +0.00%  "
$ws.Cells.Item(8, 5).Value = "  +0.51%  "
$ws.Cells.Item(9, 5).Value = "  +4.54%  "
$ws.Cells.Item(10, 5).Value = "  -1.67%  "
$ws.Cells.Item(11, 5).Value = "  +1.79%  "
$ws.Cells.Item(12, 5).Value = "  -1.63%  "
$ws.Cells.Item(13, 5).Value = "  +0.16%  "
$ws.Cells.Item(14, 5).Value = "  -4.73%  "
$ws.Cells.Item(15, 5).Value = "  -1.73%  "
$ws.Cells.Item(16, 5).Value = "  -0.85%  "
$ws.Cells.Item(17, 5).Value = "  -1.97%  "
$ws.Cells.Item(18, 5).Value = "  +0.23%  "
$ws.Cells.Item(19, 5).Value = "  -0.34%  "
$ws.Cells.Item(20, 5).Value = "  -4.72%  "
$ws.Cells.Item(21, 5).Value = "  -2.46%  "
$ws.Cells.Item(22, 5).Value = "  +0.28%  "
$ws.Cells.Item(23, 5).Value = "  +0.20%  "
$ws.Cells.Item(24, 5).Value = "  +0.02%  "
$ws.Cells.Item(25, 5).Value = "  -1.42%  "
$ws.Cells.Item(26, 5).Value = "  -4.10%  "
$ws.Cells.Item(27, 5).Value = "  -0.64%  "
$ws.Cells.Item(28, 5).Value = "  +0.23%  "
$ws.Cells.Item(29, 5).Value = "  -0.05%  "
$ws.Cells.Item(30, 5).Value = "  -4.00%  "
$ws.Cells.Item(31, 5).Value = "  -0.96%  "
$ws.Cells.Item(32, 5).Value = "  -0.07%  "
$ws.Cells.Item(33, 5).Value = "  -2.99%  "
$ws.Cells.Item(34, 5).Value = "  -2.76%  "
$ws.Cells.Item(35, 5).Value = "  +2.48%  "
$ws.Cells.Item(36, 5).Value = "  +1.07%  "
$ws.Cells.Item(37, 5).Value = "  -0.55%  "
$ws.Cells.Item(38, 5).Value = "  -4.10%  "
$ws.Cells.Item(39, 5).Value = "  -3.96%  "
$ws.Cells.Item(40, 5).Value = "  -1.63%  "
$ws.Cells.Item(41, 5).Value = "  +1.39%  "
$ws.Cells.Item(42, 5).Value = "  +0.53%  "
$ws.Cells.Item(43, 5).Value = "  -3.10%  "
$ws.Cells.Item(44, 5).Value = "  -1.46%  "
$ws.Cells.Item(45, 5).Value = "  -3.40%  "
$ws.Cells.Item(46, 5).Value = "  -5.39%  "
$ws.Cells.Item(47, 5).Value = "  -2.98%  "
$ws.Cells.Item(48, 5).Value = "  +1.29%  "
$ws.Cells.Item(49, 5).Value = "  +0.17%  "
$ws.Cells.Item(50, 5).Value = "  -3.05%  "
$ws.Cells.Item(51, 5).Value = "  -6.16%  "
